$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -7.080899999999993
$ws.Range("C7").Value = -13.43939999999999
$ws.Range("A8").Value = -22.3347
$ws.Range("A10").Value = -21.62259999999998
$ws.Range("E10").Value = 16
$ws.Range("A12").Value = -21.52630000000001
$ws.Range("E12").Value = 18.0068
$ws.Range("E13").Value = 16.37790000000001
$ws.Range("E14").Value = 16.81250000000001
$ws.Range("C15").Value = -14.35579999999998
$ws.Range("A18").Value = -21.9439
$ws.Range("C18").Value = -12.7364
$ws.Range("D18").Value = -7.912599999999995
$ws.Range("D19").Value = -8.807799999999993
$ws.Range("C20").Value = -12.45789999999999
$ws.Range("D27").Value = -8.512400000000001
$ws.Range("C29").Value = -11.83700000000001
$ws.Range("E29").Value = 17.28940000000001
$ws.Range("C30").Value = -12.9142
$ws.Range("C31").Value = -12.52739999999999
$ws.Range("D31").Value = -9.079200000000002
$ws.Range("E32").Value = 16.122
$ws.Range("E35").Value = 16.6406
$ws.Range("A37").Value = -20.78200000000001
$ws.Range("D38").Value = -8.4308
$ws.Range("C40").Value = -13.26440000000001
$ws.Range("D42").Value = -8.748299999999995
$ws.Range("E43").Value = 17.5448
$ws.Range("D44").Value = -7.545599999999999
$ws.Range("D47").Value = -7.513100000000001
$ws.Range("E48").Value = 17.60900000000002
$ws.Range("E49").Value = 15.88999999999999
$ws.Range("C50").Value = -13.6079
$ws.Range("E50").Value = 16.51329999999999
$ws.Range("A55").Value = -22.419
$ws.Range("E56").Value = 16.2532
$ws.Range("D58").Value = -8.206399999999995
$ws.Range("D65").Value = -7.764299999999997
$ws.Range("A68").Value = -21.676
$ws.Range("C68").Value = -11.6077
$ws.Range("E69").Value = 17.65750000000002
$ws.Range("D73").Value = -7.677499999999998
$ws.Range("C76").Value = -12.2746
$ws.Range("A77").Value = -21.27469999999999
$ws.Range("A78").Value = -20.89649999999999
$ws.Range("A81").Value = -21.9277
$ws.Range("E81").Value = 16.40699999999999
$ws.Range("A82").Value = -21.93770000000001
$ws.Range("C87").Value = -13.87699999999999
$ws.Range("C88").Value = -13.23409999999999
$ws.Range("D90").Value = -8.160700000000004
$ws.Range("E92").Value = 18.38230000000002
$ws.Range("D94").Value = -6.766499999999999
$ws.Range("D95").Value = -7.8686
$ws.Range("C96").Value = -12.62310000000001
$ws.Range("C98").Value = -12.1777
$ws.Range("C101").Value = -13.61180000000001
$ws.Range("D101").Value = -8.132999999999994
$ws.Range("C102").Value = -12.94970000000001
